$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-9.28%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'40.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.32%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.043"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.07303"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-5.73%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.277"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.79%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.566"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-10.13%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9175"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-3.05%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1170"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-5.85%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1717"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-7.87%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08597"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-6.41%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04169"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.00%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1053"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.14%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001279"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-1.04%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005770"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.97%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'1.30%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "'0.3276"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-2.44%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.842"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-6.62%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1350"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.21%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2882"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.04%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.03862"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-4.16%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001267"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.26%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.003852"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-6.47%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001281"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.63%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003722"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Value = "'0.02313"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'0.04972"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-7.14%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.006721"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'237.05%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007675"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.52%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1277"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-3.15%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007358"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'4.46%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007065"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-14.77%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.2888"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-16.31%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006429"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-3.70%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.008538"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-95.69%"
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.004202"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.03%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("D51").Style = "Normal"
